$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update defect counts on the "results" sheet (reducing error counts, likely
# because "errors" are now represented with icons instead of being counted).

$ws.Range("G4").Value = 0
$ws.Range("J4").Value = 1

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("J6").Value = 5

$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("J9").Value = 2

$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("J10").Value = 4
